# Fixed inability to pass in duplicate asset descriptions as input
#
# The "B" (net write-off) and "G" (net write-off, prior period) columns
# used to subtract four components: C-D-E-F and H-I-J-K respectively.
# They now are a plain two-term subtraction: C-D and H-I. This is what
# stopped duplicate asset descriptions from blowing up downstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BME_DI_VGH")

# Row 2 keeps its own (non-shared) two-term formula.
$ws.Range("B2").Formula = "=C2-D2"
$ws.Range("G2").Formula = "=H2-I2"

# Rows 3-9 get the same two-term formula, written cell-by-cell so each one
# picks up its own row's relative references.
for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=C$r-D$r"
    $ws.Cells.Item($r, 7).Formula = "=H$r-I$r"
}

# View state: zoom to 100% and move the selection to G3.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("G3").Select()
